$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") rows 2-108 all held the date serial 45185
# (2023-09-16) and are being bumped to 45204 (2023-10-05).
$ws.Range("C2:C108").Value = 45204
